# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - Dinamarca's stats overtook India's and Chile's, so those three rows
#   are re-ranked (row 27 = Dinamarca, row 28 = India, row 29 = Chile)
# - refresh the numeric counters for several other countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 10:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 400549
$ws.Range("C4").Value = 214
$ws.Range("E4").Value = 365981

# Austria (row 18)
$ws.Range("B18").Value = 12721
$ws.Range("C18").Value = 82
$ws.Range("E18").Value = 8432

# Re-ranked block: Dinamarca moves above India and Chile
# Row 27 becomes Dinamarca
$ws.Range("A27").Value = "Dinamarca"
$ws.Range("B27").Value = 5386
$ws.Range("C27").Value = 315
$ws.Range("D27").Value = 1491
$ws.Range("E27").Value = 3692
$ws.Range("F27").Value = 127
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 203

# Row 28 becomes India
$ws.Range("A28").Value = "India"
$ws.Range("B28").Value = 5360
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 468
$ws.Range("E28").Value = 4728
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 164

# Row 29 becomes Chile
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 5116
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 898
$ws.Range("E29").Value = 4175
$ws.Range("F29").Value = 337
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 43

# Polonia (row 31)
$ws.Range("B31").Value = 5000
$ws.Range("C31").Value = 152
$ws.Range("E31").Value = 4673
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = 136

# row 37
$ws.Range("B37").Value = 3870
$ws.Range("C37").Value = 106
$ws.Range("D37").Value = 96
$ws.Range("E37").Value = 3592
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 182

# row 76
$ws.Range("B76").Value = 709
$ws.Range("C76").Value = 12
$ws.Range("D76").Value = 53
$ws.Range("E76").Value = 649
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 7

# row 105
$ws.Range("B105").Value = 263
$ws.Range("C105").Value = 2
$ws.Range("E105").Value = 220
